# Implemented Strategy pattern for file export.
# Appends the new "Strategy" export-strategy transaction row (row 5) to the
# Transactions sheet, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the five new cells for row 5. A leading "'" forces Excel's
# quote-prefix / text entry semantics for the string cells, which is how
# this workbook's existing data rows (2-4) ended up carrying cell style
# index 1 (quotePrefix="1") in styles.xml - keeps the new row visually /
# structurally consistent with the rows above it.
$ws.Range("A5").Value = "'GMWOGER643"
$ws.Range("B5").Value = 666
$ws.Range("C5").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D5").Value = "'IrregularExpense"
$ws.Range("E5").Value = "'Strategy"

# Copy the formatting (cell style) of the row above down onto the new row so
# every cell - including the numeric Amount cell, which the quote-prefix
# trick above doesn't reach - ends up on style index 1, same as row 4.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

# Recalculate so any dependent state (dimension, used range, etc.) is fresh
# before save - mirrors the workbook's fullCalcOnLoad intent.
$excel.CalculateFull()
